$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Papel" column header (E1) -> clear it (column removed from data / merged into D)
$ws.Range("E1").Value = $null

# Update D2 value: combine the two "papel"/role lists into a single list, drop "usuario"/extra values
$ws.Range("D2").Value = "[ administrador | servidor | aluno ]"

# Clear E2 (previously held the second role list)
$ws.Range("E2").Value = $null
